$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Title-case Spanish connector words (de/del/la/las/el/los/y) in state & municipality names ---
$ws.Range('B5').Value = 'Pabellón De Arteaga'
$ws.Range('B6').Value = 'Rincón De Romos'
$ws.Range('B21').Value = 'Amatenango Del Valle'
$ws.Range('B24').Value = 'Bejucal De Ocampo'
$ws.Range('B27').Value = 'Chiapa De Corzo'
$ws.Range('B30').Value = 'Comitán De Domínguez'
$ws.Range('B41').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B47').Value = 'Salto De Agua'
$ws.Range('B48').Value = 'San Cristóbal De Las Casas'
$ws.Range('B69').Value = 'Coyame Del Sotol'
$ws.Range('B75').Value = 'Hidalgo Del Parral'
$ws.Range('B101').Value = 'San Juan De Sabinas'
$ws.Range('B107').Value = 'Villa De Álvarez'
$ws.Range('A109').Value = 'Ciudad De México'
$ws.Range('B113').Value = 'Cuajimalpa De Morelos'
$ws.Range('B138').Value = 'Pánuco De Coronado'
$ws.Range('B142').Value = 'San Juan De Guadalupe'
$ws.Range('A150').Value = 'Estado De México'
$ws.Range('B150').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B152').Value = 'Almoloya De Alquisiras'
$ws.Range('B153').Value = 'Almoloya De Juárez'
$ws.Range('B156').Value = 'Atizapán De Zaragoza'
$ws.Range('B163').Value = 'Ecatepec De Morelos'
$ws.Range('B165').Value = 'Ixtapan De La Sal'
$ws.Range('B171').Value = 'Naucalpan De Juárez'
$ws.Range('B183').Value = 'Tlalnepantla De Baz'
$ws.Range('B187').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B188').Value = 'Villa De Allende'
$ws.Range('B195').Value = 'San Miguel De Allende'
$ws.Range('B196').Value = 'Apaseo El Alto'
$ws.Range('B197').Value = 'Apaseo El Grande'
$ws.Range('B203').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B206').Value = 'Jaral Del Progreso'
$ws.Range('B217').Value = 'San Diego De La Unión'
$ws.Range('B219').Value = 'San Francisco Del Rincón'
$ws.Range('B221').Value = 'San Luis De La Paz'
$ws.Range('B222').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B223').Value = 'Silao De La Victoria'
$ws.Range('B225').Value = 'Valle De Santiago'
$ws.Range('B231').Value = 'Acapulco De Juárez'
$ws.Range('B233').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B234').Value = 'Alcozauca De Guerrero'
$ws.Range('B237').Value = 'Atenango Del Río'
$ws.Range('B238').Value = 'Atoyac De Álvarez'
$ws.Range('B239').Value = 'Ayutla De Los Libres'
$ws.Range('B241').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B243').Value = 'Coyuca De Benítez'
$ws.Range('B245').Value = 'Cuetzala Del Progreso'
$ws.Range('B246').Value = 'Cutzamala De Pinzón'
$ws.Range('B250').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B261').Value = 'Taxco De Alarcón'
$ws.Range('B263').Value = 'Técpan De Galeana'
$ws.Range('B265').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B267').Value = 'Tixtla De Guerrero'
$ws.Range('B269').Value = 'Tlapa De Comonfort'
$ws.Range('B275').Value = 'Atotonilco El Grande'
$ws.Range('B280').Value = 'Huasca De Ocampo'
$ws.Range('B282').Value = 'Huejutla De Reyes'
$ws.Range('B285').Value = 'Jacala De Ledezma'
$ws.Range('B291').Value = 'Mixquiahuala De Juárez'
$ws.Range('B293').Value = 'Pachuca De Soto'
$ws.Range('B295').Value = 'Progreso De Obregón'
$ws.Range('B298').Value = 'Santiago De Anaya'
$ws.Range('B301').Value = 'Tenango De Doria'
$ws.Range('B303').Value = 'Tepehuacán De Guerrero'
$ws.Range('B304').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B307').Value = 'Tula De Allende'
$ws.Range('B308').Value = 'Tulancingo De Bravo'
$ws.Range('B309').Value = 'Villa De Tezontepec'
$ws.Range('B322').Value = 'Encarnación De Díaz'
$ws.Range('B326').Value = 'Jilotlán De Los Dolores'
$ws.Range('B332').Value = 'Ojuelos De Jalisco'
$ws.Range('B335').Value = 'Talpa De Allende'
$ws.Range('B336').Value = 'Tamazula De Gordiano'
$ws.Range('B340').Value = 'Tepatitlán De Morelos'
$ws.Range('B341').Value = 'Tizapán El Alto'
$ws.Range('B388').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B415').Value = 'Tetela Del Volcán'
$ws.Range('B416').Value = 'Tlaltizapán De Zapata'
$ws.Range('B424').Value = 'Santa María Del Oro'
$ws.Range('B439').Value = 'Mier Y Noriega'
$ws.Range('B444').Value = 'San Nicolás De Los Garza'
$ws.Range('B449').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B451').Value = 'Ayoquezco De Aldama'
$ws.Range('B454').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B455').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B459').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B460').Value = 'Oaxaca De Juárez'
$ws.Range('B461').Value = 'Ocotlán De Morelos'
$ws.Range('B462').Value = 'Pinotepa De Don Luis'
$ws.Range('B474').Value = 'San José Del Progreso'
$ws.Range('B478').Value = 'San Juan Bautista Lo De Soto'
$ws.Range('B482').Value = 'San Juan Del Río'
$ws.Range('B491').Value = 'San Martín De Los Cansecos'
$ws.Range('B492').Value = 'San Miguel El Grande'
$ws.Range('B518').Value = 'Santo Domingo De Morelos'
$ws.Range('B522').Value = 'Totontepec Villa De Morelos'
$ws.Range('B523').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B524').Value = 'Villa Sola De Vega'
$ws.Range('B525').Value = 'Zimatlán De Álvarez'
$ws.Range('B539').Value = 'Cuapiaxtla De Madero'
$ws.Range('B541').Value = 'Cuayuca De Andrade'
$ws.Range('B542').Value = 'Cuetzalan Del Progreso'
$ws.Range('B545').Value = 'Huehuetlán El Chico'
$ws.Range('B550').Value = 'Izúcar De Matamoros'
$ws.Range('B554').Value = 'Los Reyes De Juárez'
$ws.Range('B557').Value = 'Palmar De Bravo'
$ws.Range('B569').Value = 'Tepexi De Rodríguez'
$ws.Range('B571').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B574').Value = 'Tuzamapan De Galeana'
$ws.Range('B579').Value = 'Xochitlán De Vicente Suárez'
$ws.Range('B586').Value = 'Amealco De Bonfil'
$ws.Range('B588').Value = 'Cadereyta De Montes'
$ws.Range('B593').Value = 'Jalpan De Serra'
$ws.Range('B594').Value = 'Landa De Matamoros'
$ws.Range('B596').Value = 'Pinal De Amoles'
$ws.Range('B599').Value = 'San Juan Del Río'
$ws.Range('B607').Value = 'Axtla De Terrazas'
$ws.Range('B612').Value = 'Ciudad Del Maíz'
$ws.Range('B621').Value = 'Mexquitic De Carmona'
$ws.Range('B626').Value = 'San Ciro De Acosta'
$ws.Range('B632').Value = 'Santa María Del Río'
$ws.Range('B634').Value = 'Soledad De Graciano Sánchez'
$ws.Range('B639').Value = 'Tanquián De Escobedo'
$ws.Range('B643').Value = 'Villa De Arista'
$ws.Range('B644').Value = 'Villa De Arriaga'
$ws.Range('B645').Value = 'Villa De Guadalupe'
$ws.Range('B646').Value = 'Villa De Ramos'
$ws.Range('B647').Value = 'Villa De Reyes'
$ws.Range('B682').Value = 'Jalpa De Méndez'
$ws.Range('B709').Value = 'Soto La Marina'
$ws.Range('B721').Value = 'San Pablo Del Monte'
$ws.Range('B739').Value = 'Amatlán De Los Reyes'
$ws.Range('B743').Value = 'Boca Del Río'
$ws.Range('B744').Value = 'Castillo De Teayo'
$ws.Range('B746').Value = 'Cazones De Herrera'
$ws.Range('B754').Value = 'Cosamaloapan De Carpio'
$ws.Range('B755').Value = 'Cosautlán De Carvajal'
$ws.Range('B765').Value = 'Hueyapan De Ocampo'
$ws.Range('B766').Value = 'Ignacio De La Llave'
$ws.Range('B774').Value = 'Juchique De Ferrer'
$ws.Range('B777').Value = 'Lerdo De Tejada'
$ws.Range('B779').Value = 'Martínez De La Torre'
$ws.Range('B780').Value = 'Medellín De Bravo'
$ws.Range('B789').Value = 'Ozuluama De Mascareñas'
$ws.Range('B793').Value = 'Poza Rica De Hidalgo'
$ws.Range('B799').Value = 'Sayula De Alemán'
$ws.Range('B800').Value = 'Soledad De Doblado'
$ws.Range('B831').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B833').Value = 'Concepción Del Oro'
$ws.Range('B846').Value = 'Moyahua De Estrada'

# --- Remove trailing footer/metadata rows (860-864), matching new dimension A1:D858 ---
$ws.Rows('860:864').Delete()
